# Update the "想去人数" (number of people interested) figures for three
# events in both the "展览" sheet and the duplicated "全部类型" sheet.
#
# F2: 525 -> 529
# F4: 167 -> 168
# F7: 738 -> 739

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 529
    $ws.Range("F4").Value = 168
    $ws.Range("F7").Value = 739
}
